$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.096.71"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "3.935.68"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D5").Value = "485.51"
$ws.Range("D6").Value = "148.31"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.725"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("E10").Value = "  +9.46%  "
$ws.Range("D11").Value = "0.0000356"
$ws.Range("E11").Value = "  +13.30%  "
$ws.Range("D12").Value = "42.67"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "10.53"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "4.570.89"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "3.956.30"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "19.74"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").Value = "69.184.42"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "436.21"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "14.59"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "87.43"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  +13.87%  "
$ws.Range("D26").Value = "3.57"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "10.52"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "38.22"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "5.89"
$ws.Range("E29").Value = "  +6.44%  "
$ws.Range("D30").Value = "714.52"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").Value = "13.24"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("D34").Value = "0.0₃0912"
$ws.Range("E34").Value = "  +34.01%  "
$ws.Range("D35").Value = "41.37"
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("D36").Value = "58.68"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  -6.92%  "
$ws.Range("D38").Value = "5.54"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  +6.75%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E42").Value = "  +6.83%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "0.340"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "147.74"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("D50").Value = "3.14"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  -2.71%  "
